$wb = $excel.ActiveWorkbook

# Sheet ALC, row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 779.3333
$ws.Range("I19").Value = 849
$ws.Range("J19").Value = 744.5
$ws.Range("K19").Value = 849
$ws.Range("L19").Value = 744.5
$ws.Range("M19").Value = -674
$ws.Range("N19").Value = -1094.5

# Sheet ALC, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2365.7144
$ws.Range("I112").Value = 1349.5
$ws.Range("K112").Value = 4048.5
$ws.Range("M112").Value = -2940.5

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4991.567
$ws.Range("I132").Value = 4731.3706
$ws.Range("J132").Value = 7333.3335
$ws.Range("K132").Value = 14194.1118
$ws.Range("L132").Value = 22000.0005
$ws.Range("M132").Value = -11664.1118
$ws.Range("N132").Value = -27060.0005

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 11173.174
$ws.Range("I137").Value = 24924.111
$ws.Range("J137").Value = 2333.2856
$ws.Range("K137").Value = 74772.333
$ws.Range("L137").Value = 6999.8568
$ws.Range("M137").Value = -72222.333
$ws.Range("N137").Value = -12099.8568

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1368.7428
$ws.Range("I138").Value = 833.03845
$ws.Range("J138").Value = 2916.3333
$ws.Range("K138").Value = 2499.11535
$ws.Range("L138").Value = 8748.999899999999
$ws.Range("M138").Value = 2640.88465
$ws.Range("N138").Value = -19028.9999

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 52667.242
$ws.Range("J45").Value = 3854.4092
$ws.Range("L45").Value = 3854.4092
$ws.Range("N45").Value = -4608.4092

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8310.917
$ws.Range("I61").Value = 13280.546
$ws.Range("J61").Value = 4105.846
$ws.Range("K61").Value = 13280.546
$ws.Range("L61").Value = 4105.846
$ws.Range("M61").Value = -13068.546
$ws.Range("N61").Value = -4529.846

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 17828.715
$ws.Range("I74").Value = 51750
$ws.Range("K74").Value = 51750
$ws.Range("M74").Value = -50876

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 17828.715
$ws.Range("I77").Value = 51750
$ws.Range("K77").Value = 258750
$ws.Range("M77").Value = -254382

# Sheet ARM, row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2704.4167
$ws.Range("I110").Value = 1744.125
$ws.Range("K110").Value = 1744.125
$ws.Range("M110").Value = 300.875

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2006996.8
$ws.Range("I122").Value = 8277.286
$ws.Range("J122").Value = 3755876.5
$ws.Range("K122").Value = 24831.858
$ws.Range("L122").Value = 11267629.5
$ws.Range("M122").Value = -22381.858
$ws.Range("N122").Value = -11272529.5

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3751.88
$ws.Range("I132").Value = 4236.2666
$ws.Range("J132").Value = 3025.3
$ws.Range("K132").Value = 12708.7998
$ws.Range("L132").Value = 9075.900000000001
$ws.Range("M132").Value = -10178.7998
$ws.Range("N132").Value = -14135.9

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 8310.917
$ws.Range("I136").Value = 13280.546
$ws.Range("J136").Value = 4105.846
$ws.Range("K136").Value = 39841.638
$ws.Range("L136").Value = 12317.538
$ws.Range("M136").Value = -37291.638
$ws.Range("N136").Value = -17417.538

# Sheet BSM, row 98
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# Sheet CRP, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1700
$ws.Range("I16").Value = 1566.6666
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 1566.6666
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -1279.6666
$ws.Range("N16").Value = -3074

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18553.666
$ws.Range("I31").Value = 42000
$ws.Range("J31").Value = 6830.5
$ws.Range("K31").Value = 42000
$ws.Range("L31").Value = 6830.5
$ws.Range("M31").Value = -41705
$ws.Range("N31").Value = -7420.5

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 18553.666
$ws.Range("I34").Value = 42000
$ws.Range("J34").Value = 6830.5
$ws.Range("K34").Value = 42000
$ws.Range("L34").Value = 6830.5
$ws.Range("M34").Value = -41798
$ws.Range("N34").Value = -7234.5

# Sheet CRP, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5810981.5
$ws.Range("I99").Value = 10559785
$ws.Range("K99").Value = 10559785
$ws.Range("M99").Value = -10558287

# Sheet CRP, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1700
$ws.Range("I113").Value = 1566.6666
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1566.6666
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 603.3334
$ws.Range("N113").Value = -6840

# Sheet CRP, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5810981.5
$ws.Range("I126").Value = 10559785
$ws.Range("K126").Value = 31679355
$ws.Range("M126").Value = -31676885

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3057.5833
$ws.Range("I132").Value = 2854.6667
$ws.Range("K132").Value = 8564.000100000001
$ws.Range("M132").Value = -6034.000100000001

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 9253
$ws.Range("I134").Value = 16931.625
$ws.Range("K134").Value = 50794.875
$ws.Range("M134").Value = -48259.875

# Sheet CUL, row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5863.5713
$ws.Range("I56").Value = 5863.5713
$ws.Range("K56").Value = 5863.5713
$ws.Range("M56").Value = -5333.5713

# Sheet CUL, row 99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 8550
$ws.Range("J99").Value = 9908.333
$ws.Range("L99").Value = 29724.999
$ws.Range("N99").Value = -34216.999

# Sheet GSM, row 29
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 3000
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11661.417
$ws.Range("I80").Value = 17120
$ws.Range("J80").Value = 4019.4
$ws.Range("K80").Value = 17120
$ws.Range("L80").Value = 4019.4
$ws.Range("M80").Value = -16122
$ws.Range("N80").Value = -6015.4

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 11661.417
$ws.Range("I83").Value = 17120
$ws.Range("J83").Value = 4019.4
$ws.Range("K83").Value = 85600
$ws.Range("L83").Value = 20097
$ws.Range("M83").Value = -80608
$ws.Range("N83").Value = -30081

# Sheet GSM, row 105
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 89333
$ws.Range("J105").Value = 89333
$ws.Range("L105").Value = 89333
$ws.Range("N105").Value = -96321

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9943.381
$ws.Range("I122").Value = 7448
$ws.Range("K122").Value = 22344
$ws.Range("M122").Value = -19894

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4573.2563
$ws.Range("I132").Value = 5280.7144
$ws.Range("J132").Value = 2772.4546
$ws.Range("K132").Value = 15842.1432
$ws.Range("L132").Value = 8317.3638
$ws.Range("M132").Value = -13312.1432
$ws.Range("N132").Value = -13377.3638

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 996543.44
$ws.Range("I132").Value = 1491790.8
$ws.Range("K132").Value = 4475372.4
$ws.Range("M132").Value = -4472842.4

# Sheet WVR, row 28
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

# Sheet WVR, row 105
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 41499.5
$ws.Range("J105").Value = 41499.5
$ws.Range("L105").Value = 41499.5
$ws.Range("N105").Value = -48487.5

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4187.0835
$ws.Range("I122").Value = 2086.5293
$ws.Range("K122").Value = 6259.5879
$ws.Range("M122").Value = -3809.5879

# Sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 53072.25
$ws.Range("I126").Value = 134897.33
$ws.Range("K126").Value = 404691.99
$ws.Range("M126").Value = -402221.99
